$wb = $excel.ActiveWorkbook
$wsAbout = $wb.Worksheets.Item("About")
$wsData = $wb.Worksheets.Item("BBPPRTY")

# --- Data edits ---
# Turn on (suppress retrofit availability for) the CCS technologies called out
# in the commit message: lignite w/ CCS, hard coal w/ CCS, and biomass w/ CCS.
# Each is a per-year boolean flag (row) spanning years 2021-2050 (cols B:AE).
# Row 19 = "hard coal w CCS"
# Row 21 = "biomass w CCS"
# Row 22 = "lignite w CCS"
$wsData.Range("B19:AE19").Value = 1
$wsData.Range("B21:AE21").Value = 1
$wsData.Range("B22:AE22").Value = 1

# --- Selection / view state updates (match the saved cursor positions) ---
$wsAbout.Range("B20").Select() | Out-Null

$wsData.Activate() | Out-Null
$wsData.Range("B19:AE22").Select() | Out-Null
